$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 row(s) after original row 35
$ws.Range("A36:A37").EntireRow.Insert()
$ws.Range("A36").Value = "Stowaway"
$ws.Range("B36").Value = "Transport – stowaway"
$ws.Range("A37").Value = "Stowaway (L)"
$ws.Range("B37").Value = "Transport – stowaway"

# Insert 1 row(s) after original row 20
$ws.Range("A21").EntireRow.Insert()
$ws.Range("A21").Value = "Research"
$ws.Range("B21").Value = "Escape from confinement: Research and ex-situ breeding"

# Insert 3 row(s) after original row 19
$ws.Range("A20:A22").EntireRow.Insert()
$ws.Range("A20").Value = "Pet (L)"
$ws.Range("B20").Value = "Escape from confinement: Pet/aquarium/terrarium species"
$ws.Range("A21").Value = "Pet"
$ws.Range("B21").Value = "Escape from confinement: Pet/aquarium/terrarium species"
$ws.Range("A22").Value = "Ornamental Purposes"
$ws.Range("B22").Value = "Escape from confinement: Ornamental purpose other than horticulture"

# Insert 1 row(s) after original row 18
$ws.Range("A19").EntireRow.Insert()
$ws.Range("A19").Value = "Zoo"
$ws.Range("B19").Value = "Escape from confinement: Botanical garden/zoo/aquaria"

# Insert 2 row(s) after original row 15
$ws.Range("A16:A17").EntireRow.Insert()
$ws.Range("A16").Value = "Fur Farming"
$ws.Range("B16").Value = "Escape from confinement: Fur farms"
$ws.Range("A17").Value = "Wild Fur"
$ws.Range("B17").Value = "Escape from confinement: Fur farms"

# Insert 1 row(s) after original row 13
$ws.Range("A14").EntireRow.Insert()
$ws.Range("A14").Value = "Farming"
$ws.Range("B14").Value = "Escape from confinement: Farmed animals"

# Insert 1 row(s) after original row 10
$ws.Range("A11").EntireRow.Insert()
$ws.Range("A11").Value = "Fauna Improvement, Hunting"
$ws.Range("B11").Value = "Release in nature: Other intentional release"

# Insert 3 row(s) after original row 8
$ws.Range("A9:A11").EntireRow.Insert()
$ws.Range("A9").Value = "Conservation, Zoo"
$ws.Range("B9").Value = "Release in nature: Introduction for conservation purposes"
$ws.Range("A10").Value = "Conservation (L)"
$ws.Range("B10").Value = "Release in nature: Introduction for conservation purposes"
$ws.Range("A11").Value = "Conservation"
$ws.Range("B11").Value = "Release in nature: Introduction for conservation purposes"

# Insert 2 row(s) after original row 7
$ws.Range("A8:A9").EntireRow.Insert()
$ws.Range("A8").Value = "Hunting (L)"
$ws.Range("B8").Value = "Release in nature: Hunting in the wild"
$ws.Range("A9").Value = "Hunting"
$ws.Range("B9").Value = "Release in nature: Hunting in the wild"

# Insert 2 row(s) after original row 3
$ws.Range("A4:A5").EntireRow.Insert()
$ws.Range("A4").Value = "Biological Control (L)"
$ws.Range("B4").Value = "Release in nature: Biological control"
$ws.Range("A5").Value = "Biological Control"
$ws.Range("B5").Value = "Release in nature: Biological control"

# Match final selection shown in the authored workbook
$ws.Range("A26").Select()
